$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (player name, position, team) for rows 2-17
$data = @(
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Kyle Kuzma", "PF", "Milwaukee Bucks"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Bam Adebayo", "PF,C", "Miami Heat"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Nikola Jovic", "PF,C", "Miami Heat"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Rui Hachimura", "SF,PF", "Los Angeles Lakers"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Zach LaVine", "SG,SF", "Sacramento Kings"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the now-unused rows 18 and 19 (table shrank from 18 players to 16)
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(18).Delete()
